$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.627428412437439
$ws.Range("B1").Value = 0.9242424964904785
$ws.Range("C1").Value = 4.039738655090332
$ws.Range("D1").Value = 2.105851411819458
$ws.Range("E1").Value = 1.65105926990509
